{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Applies three changes to docs/controlValidation.docx:\n//  1. \"Must be in format \"wXXXXXXX\"\" \u2014 collapse the 3 runs (with the\n//     spell-check proofErr wrapper around \"wXXXXXXX\") into a single run\n//     with the same combined text.\n//  2. \"Integers >= 0 only\" -> \"Float  >= 0 only\" (two runs: \"Float \" and\n//     \" >= 0 only\").\n//  3. \"Max 2 characters\" -> \"Max 4 characters (XX.XX)\" (four runs: \"Max \",\n//     \"4\", \" characters\", \" (XX.XX)\").\n//\n// Small flat-OPC snippets are fed through Range.insertOoxml(...,\"Replace\")\n// so the exact run boundaries in the target OOXML are reproduced (plain\n// Range/Paragraph.insertText merges adjacent text into a single run, which\n// would not match runs 2 & 3 of the diff).\n\nconst body = context.document.body;\n\n// --- Change 1: merge \"Must be in format \"wXXXXXXX\"\" into one run -------\nconst wFormatHits = body.search(\"Must be in format\", { matchCase: true });\nwFormatHits.load(\"items/text\");\nawait context.sync();\n\nfor (const hit of wFormatHits.items) {\n  const para = hit.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n  if (/^Must be in format [\\u201c\"]wXXXXXXX[\\u201d\"]$/.test(para.text)) {\n    // A plain Replace naturally coalesces into a single run and drops the\n    // proofErr spell-check markers that bracketed the middle run.\n    para.insertText(\"Must be in format \\u201cwXXXXXXX\\u201d\", \"Replace\");\n  }\n}\nawait context.sync();\n\n// --- Change 2: \"Integers >= 0 only\" -> \"Float \" + \" >= 0 only\" ---------\nconst floatOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n</Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\">Float </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> &gt;= 0 only</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nconst integersHits = body.search(\"Integers >= 0 only\", { matchCase: true });\nintegersHits.load(\"items\");\nawait context.sync();\nfor (const hit of integersHits.items) {\n  hit.insertOoxml(floatOoxml, \"Replace\");\n}\nawait context.sync();\n\n// --- Change 3: \"Max 2 characters\" -> \"Max \"+\"4\"+\" characters\"+\" (XX.XX)\" -\nconst maxOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n<pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n</Relationships>\n</pkg:xmlData>\n</pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\">Max </w:t></w:r>\n<w:r><w:t>4</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> characters</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> (XX.XX)</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nconst maxHits = body.search(\"Max 2 characters\", { matchCase: true });\nmaxHits.load(\"items\");\nawait context.sync();\nfor (const hit of maxHits.items) {\n  hit.insertOoxml(maxOoxml, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Applies three changes to docs/controlValidation.docx:\n#  1. \"Must be in format \"wXXXXXXX\"\" -> collapse the 3 runs (with the\n#     spell-check proofErr wrapper around \"wXXXXXXX\") into a single run\n#     with the same combined text.\n#  2. \"Integers >= 0 only\" -> \"Float  >= 0 only\" (two runs: \"Float \" and\n#     \" >= 0 only\").\n#  3. \"Max 2 characters\" -> \"Max 4 characters (XX.XX)\" (four runs: \"Max \",\n#     \"4\", \" characters\", \" (XX.XX)\").\n#\n# $word / $d are pre-seeded by the harness ($d = $word.ActiveDocument).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: merge the \"Must be in format \"wXXXXXXX\"\" runs into one ---\n# Find/Replace with identical Find & Replacement text rewrites the matched\n# range as a single run and drops the proofErr spell-check bookmarks that\n# bracketed the middle run.\n$needle = \"Must be in format \" + [char]0x201C + \"wXXXXXXX\" + [char]0x201D\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)\n\n# --- Change 2: \"Integers >= 0 only\" -> \"Float \" + \" >= 0 only\" (2 runs) ---\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq (\"Integers >= 0 only\" + [char]13)) {\n        $r = $p.Range\n        # Exclude the trailing paragraph mark from the replaced span.\n        $sub = $d.Range($r.Start, $r.End - 1)\n        $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n               '<w:r><w:t xml:space=\"preserve\">Float </w:t></w:r>' +\n               '<w:r><w:t xml:space=\"preserve\"> &gt;= 0 only</w:t></w:r>' +\n               '</w:p>'\n        $sub.InsertXML($xml)\n        break\n    }\n}\n\n# --- Change 3: \"Max 2 characters\" -> \"Max \"+\"4\"+\" characters\"+\" (XX.XX)\" (4 runs) ---\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq (\"Max 2 characters\" + [char]13)) {\n        $r = $p.Range\n        $sub = $d.Range($r.Start, $r.End - 1)\n        $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n               '<w:r><w:t xml:space=\"preserve\">Max </w:t></w:r>' +\n               '<w:r><w:t>4</w:t></w:r>' +\n               '<w:r><w:t xml:space=\"preserve\"> characters</w:t></w:r>' +\n               '<w:r><w:t xml:space=\"preserve\"> (XX.XX)</w:t></w:r>' +\n               '</w:p>'\n        $sub.InsertXML($xml)\n        break\n    }\n}\n"}
